# Insert a new price-record row into the "Ají" sheet at row 263, pushing the
# existing rows 263-327 down to 264-328 (their original content stays intact).
# Populate the newly inserted row with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 263 (shifts rows 263..327 down to 264..328)
$ws.Rows.Item(263).Insert()

# Fill in the new row 263 with the new record's data
$ws.Cells.Item(263, 1).Value  = 9
$ws.Cells.Item(263, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(263, 3).Value  = "Metropolitana"
$ws.Cells.Item(263, 4).Value  = 44798
$ws.Cells.Item(263, 5).Value  = 13
$ws.Cells.Item(263, 6).Value  = 100112021
$ws.Cells.Item(263, 7).Value  = "Ají"
$ws.Cells.Item(263, 8).Value  = "Inferno"
$ws.Cells.Item(263, 9).Value  = "Primera"
$ws.Cells.Item(263, 10).Value = 153
$ws.Cells.Item(263, 11).Value = 13000
$ws.Cells.Item(263, 12).Value = 15000
$ws.Cells.Item(263, 13).Value = 13889
$ws.Cells.Item(263, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(263, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(263, 16).Value = 1389
$ws.Cells.Item(263, 17).Value = 10
$ws.Cells.Item(263, 18).Value = "Hortaliza"
